# ============================================================================
# feat: add 2022-Q3 data
#
# 1. Summary sheet ("总计"): insert a new data row at row 2 for "2022-Q3"
#    (C=21 holdings, D=2.37 market value), pushing the existing quarter rows
#    down by one. The running index in column A stays anchored to the row
#    position (0,1,2,...) rather than following the data.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计",
#    populated with the per-fund holding breakdown (21 funds).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Summary sheet ("总计"): shift existing rows down and insert 2022-Q3 row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push rows 2..8 down to 3..9, leaving row 2 blank (values + formatting of
# the rows below are carried along automatically).
$summary.Rows.Item(2).Insert()

# The new blank row lost any sensible formatting recipe - reset B:D to the
# sheet's default (unstyled) look, matching the other data rows.
$summary.Range("B2:D2").Style = "Normal"

# Column A carries the bold/centered/bordered "index" style (style index 2
# in the original file) on every data row - copy that formatting from the
# row below onto the new A2 cell.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New 2022-Q3 data row
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 21
$summary.Range("D2").Value = 2.37

# The running index in column A is anchored to the row position, not to the
# data that moved down with the insert - renumber rows 3..9 back to 1..7.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7

# ---------------------------------------------------------------------------
# 2) New worksheet "2022-Q3" (inserted right after "总计")
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Add($null, $summary)
$ws.Name = "2022-Q3"

# Header row (bold, centered, thin border - matches the other quarter sheets)
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"
$ws.Range("B1:H1").Font.Bold = $true
$ws.Range("B1:H1").HorizontalAlignment = -4108
$ws.Range("B1:H1").VerticalAlignment = -4160
$ws.Range("B1:H1").Borders.LineStyle = 1

# Columns B, D, E, F, G hold text-formatted numbers (leading zeros / trailing
# zeros must be preserved verbatim), so force text format before writing.
$ws.Range("B2:B22").NumberFormat = "@"
$ws.Range("D2:G22").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "340007"
$ws.Range("C2").Value = "兴全社会责任混合"
$ws.Range("D2").Value = "44.40"
$ws.Range("E2").Value = "87.93"
$ws.Range("F2").Value = "3.48"
$ws.Range("G2").Value = "1.5451"
$ws.Range("H2").Value = 10
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "009564"
$ws.Range("C3").Value = "汇安消费龙头混合A"
$ws.Range("D3").Value = "7.60"
$ws.Range("E3").Value = "93.84"
$ws.Range("F3").Value = "3.54"
$ws.Range("G3").Value = "0.2690"
$ws.Range("H3").Value = 8
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "005274"
$ws.Range("C4").Value = "中银景福回报混合A"
$ws.Range("D4").Value = "10.39"
$ws.Range("E4").Value = "23.39"
$ws.Range("F4").Value = "1.05"
$ws.Range("G4").Value = "0.1091"
$ws.Range("H4").Value = 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "006648"
$ws.Range("C5").Value = "汇安多因子混合A"
$ws.Range("D5").Value = "3.69"
$ws.Range("E5").Value = "93.49"
$ws.Range("F5").Value = "2.59"
$ws.Range("G5").Value = "0.0956"
$ws.Range("H5").Value = 9
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "006649"
$ws.Range("C6").Value = "汇安多因子混合C"
$ws.Range("D6").Value = "1.95"
$ws.Range("E6").Value = "93.49"
$ws.Range("F6").Value = "2.59"
$ws.Range("G6").Value = "0.0505"
$ws.Range("H6").Value = 9
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "007318"
$ws.Range("C7").Value = "中银民丰回报混合"
$ws.Range("D7").Value = "4.16"
$ws.Range("E7").Value = "22.58"
$ws.Range("F7").Value = "1.05"
$ws.Range("G7").Value = "0.0437"
$ws.Range("H7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "163823"
$ws.Range("C8").Value = "中银稳健策略灵活配置混合"
$ws.Range("D8").Value = "2.06"
$ws.Range("E8").Value = "45.80"
$ws.Range("F8").Value = "2.12"
$ws.Range("G8").Value = "0.0437"
$ws.Range("H8").Value = 5
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "010558"
$ws.Range("C9").Value = "汇安鑫利优选混合A"
$ws.Range("D9").Value = "1.30"
$ws.Range("E9").Value = "93.99"
$ws.Range("F9").Value = "2.60"
$ws.Range("G9").Value = "0.0338"
$ws.Range("H9").Value = 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "002535"
$ws.Range("C10").Value = "中银鑫利灵活配置混合A"
$ws.Range("D10").Value = "2.81"
$ws.Range("E10").Value = "23.75"
$ws.Range("F10").Value = "1.05"
$ws.Range("G10").Value = "0.0295"
$ws.Range("H10").Value = 5
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "008773"
$ws.Range("C11").Value = "中银景泰回报混合"
$ws.Range("D11").Value = "2.34"
$ws.Range("E11").Value = "22.68"
$ws.Range("F11").Value = "1.06"
$ws.Range("G11").Value = "0.0248"
$ws.Range("H11").Value = 5
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "006952"
$ws.Range("C12").Value = "中银景元回报混合"
$ws.Range("D12").Value = "2.19"
$ws.Range("E12").Value = "26.13"
$ws.Range("F12").Value = "1.05"
$ws.Range("G12").Value = "0.0230"
$ws.Range("H12").Value = 6
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "011858"
$ws.Range("C13").Value = "安信消费升级一年持有期混合A"
$ws.Range("D13").Value = "0.86"
$ws.Range("E13").Value = "76.21"
$ws.Range("F13").Value = "2.65"
$ws.Range("G13").Value = "0.0228"
$ws.Range("H13").Value = 8
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "002288"
$ws.Range("C14").Value = "中银稳进策略灵活配置混合A"
$ws.Range("D14").Value = "0.62"
$ws.Range("E14").Value = "64.79"
$ws.Range("F14").Value = "3.18"
$ws.Range("G14").Value = "0.0197"
$ws.Range("H14").Value = 4
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "009565"
$ws.Range("C15").Value = "汇安消费龙头混合C"
$ws.Range("D15").Value = "0.49"
$ws.Range("E15").Value = "93.84"
$ws.Range("F15").Value = "3.54"
$ws.Range("G15").Value = "0.0173"
$ws.Range("H15").Value = 8
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "003889"
$ws.Range("C16").Value = "汇安丰泽灵活配置混合A"
$ws.Range("D16").Value = "0.64"
$ws.Range("E16").Value = "93.72"
$ws.Range("F16").Value = "2.65"
$ws.Range("G16").Value = "0.0170"
$ws.Range("H16").Value = 8
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "010559"
$ws.Range("C17").Value = "汇安鑫利优选混合C"
$ws.Range("D17").Value = "0.65"
$ws.Range("E17").Value = "93.99"
$ws.Range("F17").Value = "2.60"
$ws.Range("G17").Value = "0.0169"
$ws.Range("H17").Value = 10
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "003890"
$ws.Range("C18").Value = "汇安丰泽灵活配置混合C"
$ws.Range("D18").Value = "0.28"
$ws.Range("E18").Value = "93.72"
$ws.Range("F18").Value = "2.65"
$ws.Range("G18").Value = "0.0074"
$ws.Range("H18").Value = 8
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "011859"
$ws.Range("C19").Value = "安信消费升级一年持有期混合C"
$ws.Range("D19").Value = "0.13"
$ws.Range("E19").Value = "76.21"
$ws.Range("F19").Value = "2.65"
$ws.Range("G19").Value = "0.0034"
$ws.Range("H19").Value = 8
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "002536"
$ws.Range("C20").Value = "中银鑫利灵活配置混合C"
$ws.Range("D20").Value = "0.14"
$ws.Range("E20").Value = "23.75"
$ws.Range("F20").Value = "1.05"
$ws.Range("G20").Value = "0.0015"
$ws.Range("H20").Value = 5
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "015089"
$ws.Range("C21").Value = "中银景福回报混合C"
$ws.Range("D21").Value = "0.00"
$ws.Range("E21").Value = "23.39"
$ws.Range("F21").Value = "1.05"
$ws.Range("G21").NumberFormat = "General"
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 5
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "016520"
$ws.Range("C22").Value = "中银稳进策略灵活配置混合C"
$ws.Range("D22").Value = "0.00"
$ws.Range("E22").Value = "64.79"
$ws.Range("F22").Value = "3.18"
$ws.Range("G22").NumberFormat = "General"
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 4

# Column A ("基金" running index) carries the same bold/centered/bordered
# style as the summary sheet's index column.
$ws.Range("A2:A22").Font.Bold = $true
$ws.Range("A2:A22").HorizontalAlignment = -4108
$ws.Range("A2:A22").VerticalAlignment = -4160
$ws.Range("A2:A22").Borders.LineStyle = 1
